$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '72.121.68'
$c.ClearFormats()

$c = $ws.Range('E2')
$c.NumberFormat = '@'
$c.Value = '  -0.28%  '
$c.ClearFormats()

$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '2.660.20'
$c.ClearFormats()

$c = $ws.Range('E3')
$c.NumberFormat = '@'
$c.Value = '  +1.32%  '
$c.ClearFormats()

$c = $ws.Range('E4')
$c.NumberFormat = '@'
$c.Value = '  +0.00%  '
$c.ClearFormats()

$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '595.69'
$c.ClearFormats()

$c = $ws.Range('E5')
$c.NumberFormat = '@'
$c.Value = '  -1.29%  '
$c.ClearFormats()

$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '175.29'
$c.ClearFormats()

$c = $ws.Range('E6')
$c.NumberFormat = '@'
$c.Value = '  -1.72%  '
$c.ClearFormats()

$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.523'
$c.ClearFormats()

$c = $ws.Range('E8')
$c.NumberFormat = '@'
$c.Value = '  -0.48%  '
$c.ClearFormats()

$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '2.660.20'
$c.ClearFormats()

$c = $ws.Range('E9')
$c.NumberFormat = '@'
$c.Value = '  +1.32%  '
$c.ClearFormats()

$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.169'
$c.ClearFormats()

$c = $ws.Range('E10')
$c.NumberFormat = '@'
$c.Value = '  -3.47%  '
$c.ClearFormats()

$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.170'
$c.ClearFormats()

$c = $ws.Range('E11')
$c.NumberFormat = '@'
$c.Value = '  +1.98%  '
$c.ClearFormats()

$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '0.357'
$c.ClearFormats()

$c = $ws.Range('E12')
$c.NumberFormat = '@'
$c.Value = '  -0.61%  '
$c.ClearFormats()

$c = $ws.Range('E13')
$c.NumberFormat = '@'
$c.Value = '  -0.48%  '
$c.ClearFormats()

$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '3.148.08'
$c.ClearFormats()

$c = $ws.Range('E14')
$c.NumberFormat = '@'
$c.Value = '  +1.12%  '
$c.ClearFormats()

$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '72.019.55'
$c.ClearFormats()

$c = $ws.Range('E15')
$c.NumberFormat = '@'
$c.Value = '  -0.45%  '
$c.ClearFormats()

$c = $ws.Range('E16')
$c.NumberFormat = '@'
$c.Value = '  -3.79%  '
$c.ClearFormats()

$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '26.17'
$c.ClearFormats()

$c = $ws.Range('E17')
$c.NumberFormat = '@'
$c.Value = '  -1.87%  '
$c.ClearFormats()

$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '2.655.65'
$c.ClearFormats()

$c = $ws.Range('E18')
$c.NumberFormat = '@'
$c.Value = '  +1.12%  '
$c.ClearFormats()

$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '12.38'
$c.ClearFormats()

$c = $ws.Range('E19')
$c.NumberFormat = '@'
$c.Value = '  +5.28%  '
$c.ClearFormats()

$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '7.94'
$c.ClearFormats()

$c = $ws.Range('E20')
$c.NumberFormat = '@'
$c.Value = '  -0.05%  '
$c.ClearFormats()

$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '367.20'
$c.ClearFormats()

$c = $ws.Range('E21')
$c.NumberFormat = '@'
$c.Value = '  -4.35%  '
$c.ClearFormats()

$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '4.21'
$c.ClearFormats()

$c = $ws.Range('E22')
$c.NumberFormat = '@'
$c.Value = '  +0.65%  '
$c.ClearFormats()

$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '2.06'
$c.ClearFormats()

$c = $ws.Range('E23')
$c.NumberFormat = '@'
$c.Value = '  -1.41%  '
$c.ClearFormats()

$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '71.66'
$c.ClearFormats()

$c = $ws.Range('E24')
$c.NumberFormat = '@'
$c.Value = '  -3.11%  '
$c.ClearFormats()

$c = $ws.Range('E25')
$c.NumberFormat = '@'
$c.Value = '  +0.10%  '
$c.ClearFormats()

$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '4.31'
$c.ClearFormats()

$c = $ws.Range('E26')
$c.NumberFormat = '@'
$c.Value = '  -1.89%  '
$c.ClearFormats()

$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '9.74'
$c.ClearFormats()

$c = $ws.Range('E27')
$c.NumberFormat = '@'
$c.Value = '  -3.42%  '
$c.ClearFormats()

$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '2.797.18'
$c.ClearFormats()

$c = $ws.Range('E28')
$c.NumberFormat = '@'
$c.Value = '  +1.34%  '
$c.ClearFormats()

$c = $ws.Range('E29')
$c.NumberFormat = '@'
$c.Value = '  +0.17%  '
$c.ClearFormats()

$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '0.0₃0965'
$c.ClearFormats()

$c = $ws.Range('E30')
$c.NumberFormat = '@'
$c.Value = '  +0.31%  '
$c.ClearFormats()

$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '8.18'
$c.ClearFormats()

$c = $ws.Range('E31')
$c.NumberFormat = '@'
$c.Value = '  +0.61%  '
$c.ClearFormats()

$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '496.78'
$c.ClearFormats()

$c = $ws.Range('E32')
$c.NumberFormat = '@'
$c.Value = '  -3.35%  '
$c.ClearFormats()

$c = $ws.Range('E33')
$c.NumberFormat = '@'
$c.Value = '  -2.09%  '
$c.ClearFormats()

$c = $ws.Range('E34')
$c.NumberFormat = '@'
$c.Value = '  -0.71%  '
$c.ClearFormats()

$c = $ws.Range('E35')
$c.NumberFormat = '@'
$c.Value = '  +0.02%  '
$c.ClearFormats()

$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '162.63'
$c.ClearFormats()

$c = $ws.Range('E36')
$c.NumberFormat = '@'
$c.Value = '  -0.49%  '
$c.ClearFormats()

$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '19.41'
$c.ClearFormats()

$c = $ws.Range('E37')
$c.NumberFormat = '@'
$c.Value = '  +0.53%  '
$c.ClearFormats()

$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '0.115'
$c.ClearFormats()

$c = $ws.Range('E38')
$c.NumberFormat = '@'
$c.Value = '  +2.22%  '
$c.ClearFormats()

$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '18.89'
$c.ClearFormats()

$c = $ws.Range('E39')
$c.NumberFormat = '@'
$c.Value = '  -0.95%  '
$c.ClearFormats()

$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '1.37'
$c.ClearFormats()

$c = $ws.Range('E40')
$c.NumberFormat = '@'
$c.Value = '  -2.79%  '
$c.ClearFormats()

$c = $ws.Range('B41')
$c.NumberFormat = '@'
$c.Value = 'USDe'
$c.ClearFormats()

$c = $ws.Range('C41')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$c.ClearFormats()

$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.ClearFormats()

$c = $ws.Range('E41')
$c.NumberFormat = '@'
$c.Value = '  -0.04%  '
$c.ClearFormats()

$c = $ws.Range('B42')
$c.NumberFormat = '@'
$c.Value = 'Stacks'
$c.ClearFormats()

$c = $ws.Range('C42')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$c.ClearFormats()

$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '1.74'
$c.ClearFormats()

$c = $ws.Range('E42')
$c.NumberFormat = '@'
$c.Value = '  -5.54%  '
$c.ClearFormats()

$c = $ws.Range('E43')
$c.NumberFormat = '@'
$c.Value = '  -0.38%  '
$c.ClearFormats()

$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '4.97'
$c.ClearFormats()

$c = $ws.Range('E44')
$c.NumberFormat = '@'
$c.Value = '  -2.01%  '
$c.ClearFormats()

$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '0.331'
$c.ClearFormats()

$c = $ws.Range('E45')
$c.NumberFormat = '@'
$c.Value = '  -0.62%  '
$c.ClearFormats()

$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '156.96'
$c.ClearFormats()

$c = $ws.Range('E46')
$c.NumberFormat = '@'
$c.Value = '  +4.64%  '
$c.ClearFormats()

$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '39.15'
$c.ClearFormats()

$c = $ws.Range('E47')
$c.NumberFormat = '@'
$c.Value = '  -0.67%  '
$c.ClearFormats()

$c = $ws.Range('E48')
$c.NumberFormat = '@'
$c.Value = '  +2.48%  '
$c.ClearFormats()

$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '3.71'
$c.ClearFormats()

$c = $ws.Range('E49')
$c.NumberFormat = '@'
$c.Value = '  +0.11%  '
$c.ClearFormats()

$c = $ws.Range('E50')
$c.NumberFormat = '@'
$c.Value = '  +1.72%  '
$c.ClearFormats()

$c = $ws.Range('B51')
$c.NumberFormat = '@'
$c.Value = 'Cronos'
$c.ClearFormats()

$c = $ws.Range('C51')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$c.ClearFormats()

$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '0.0756'
$c.ClearFormats()

$c = $ws.Range('E51')
$c.NumberFormat = '@'
$c.Value = '  -1.25%  '
$c.ClearFormats()
